$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Temperature (°C)"
$ws.Range("B2").Value = 22.51764705882353

$ws.Range("A3").Value = "Duration (min)"
$ws.Range("B3").Value = 53.55330882352941

$ws.Range("A4").Value = "Temperature (°C)*Duration (min)"
$ws.Range("B4").Value = 23.02389705882356

$ws.Range("B5").Value = 48.28235294117639
